$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell (outside the data table) used to stage each new value as
# plain Text. Writing text-looking numbers (e.g. "142.38", "34.00",
# "0.999") straight into Range.Value would let Excel auto-convert them
# to real numbers, losing the exact string formatting the source data
# uses (trailing zeros, dot thousand-separators, etc). Staging the value
# in a Text-formatted cell and then Copy / PasteSpecial(values-only) into
# the real target cell carries over the literal text without also
# stamping the target cell with the scratch cell's Text number format.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$scratch.Value = "61.123.40"
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)

$scratch.Value = "  -1.82%  "
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)

$scratch.Value = "2.976.60"
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)

$scratch.Value = "  -0.50%  "
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Value = "0.999"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)

$scratch.Value = "  +0.05%  "
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Value = "595.37"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)

$scratch.Value = "  +2.48%  "
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.Value = "142.38"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)

$scratch.Value = "  -2.40%  "
$scratch.Copy()
$ws.Range("E6").PasteSpecial(-4163)

$scratch.Value = "  +0.25%  "
$scratch.Copy()
$ws.Range("E7").PasteSpecial(-4163)

$scratch.Value = "2.975.87"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)

$scratch.Value = "  -0.51%  "
$scratch.Copy()
$ws.Range("E9").PasteSpecial(-4163)

$scratch.Value = "Dogecoin"
$scratch.Copy()
$ws.Range("B10").PasteSpecial(-4163)

$scratch.Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$scratch.Copy()
$ws.Range("C10").PasteSpecial(-4163)

$scratch.Value = "0.145"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)

$scratch.Value = "  -1.89%  "
$scratch.Copy()
$ws.Range("E10").PasteSpecial(-4163)

$scratch.Value = "Toncoin"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)

$scratch.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$scratch.Copy()
$ws.Range("C11").PasteSpecial(-4163)

$scratch.Value = "6.00"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)

$scratch.Value = "  +6.37%  "
$scratch.Copy()
$ws.Range("E11").PasteSpecial(-4163)

$scratch.Value = "0.454"
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)

$scratch.Value = "  +2.90%  "
$scratch.Copy()
$ws.Range("E12").PasteSpecial(-4163)

$scratch.Value = "0.0000226"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)

$scratch.Value = "  -0.48%  "
$scratch.Copy()
$ws.Range("E13").PasteSpecial(-4163)

$scratch.Value = "34.00"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)

$scratch.Value = "  -1.51%  "
$scratch.Copy()
$ws.Range("E14").PasteSpecial(-4163)

$scratch.Value = "0.125"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)

$scratch.Value = "  +2.18%  "
$scratch.Copy()
$ws.Range("E15").PasteSpecial(-4163)

$scratch.Value = "3.465.70"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)

$scratch.Value = "  -0.41%  "
$scratch.Copy()
$ws.Range("E16").PasteSpecial(-4163)

$scratch.Value = "61.089.80"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)

$scratch.Value = "  -1.84%  "
$scratch.Copy()
$ws.Range("E17").PasteSpecial(-4163)

$scratch.Value = "6.86"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)

$scratch.Value = "  -2.37%  "
$scratch.Copy()
$ws.Range("E18").PasteSpecial(-4163)

$scratch.Value = "2.972.46"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)

$scratch.Value = "447.19"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)

$scratch.Value = "  -2.05%  "
$scratch.Copy()
$ws.Range("E20").PasteSpecial(-4163)

$scratch.Value = "14.08"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)

$scratch.Value = "  +1.65%  "
$scratch.Copy()
$ws.Range("E21").PasteSpecial(-4163)

$scratch.Value = "0.679"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)

$scratch.Value = "  +0.19%  "
$scratch.Copy()
$ws.Range("E22").PasteSpecial(-4163)

$scratch.Value = "7.27"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)

$scratch.Value = "  -0.19%  "
$scratch.Copy()
$ws.Range("E23").PasteSpecial(-4163)

$scratch.Value = "82.00"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$scratch.Value = "  +2.70%  "
$scratch.Copy()
$ws.Range("E24").PasteSpecial(-4163)

$scratch.Value = "  -5.16%  "
$scratch.Copy()
$ws.Range("E25").PasteSpecial(-4163)

$scratch.Value = "RenderToken"
$scratch.Copy()
$ws.Range("B26").PasteSpecial(-4163)

$scratch.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$scratch.Copy()
$ws.Range("C26").PasteSpecial(-4163)

$scratch.Value = "10.33"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)

$scratch.Value = "  +3.54%  "
$scratch.Copy()
$ws.Range("E26").PasteSpecial(-4163)

$scratch.Value = "InternetComputer(DFINITY)"
$scratch.Copy()
$ws.Range("B27").PasteSpecial(-4163)

$scratch.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$scratch.Copy()
$ws.Range("C27").PasteSpecial(-4163)

$scratch.Value = "11.89"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)

$scratch.Value = "  -2.75%  "
$scratch.Copy()
$ws.Range("E27").PasteSpecial(-4163)

$scratch.Value = "  +0.00%  "
$scratch.Copy()
$ws.Range("E28").PasteSpecial(-4163)

$scratch.Value = "  +2.68%  "
$scratch.Copy()
$ws.Range("E29").PasteSpecial(-4163)

$scratch.Value = "0.999"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)

$scratch.Value = "  +0.00%  "
$scratch.Copy()
$ws.Range("E30").PasteSpecial(-4163)

$scratch.Value = "7.08"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)

$scratch.Value = "  -2.29%  "
$scratch.Copy()
$ws.Range("E31").PasteSpecial(-4163)

$scratch.Value = "  -2.40%  "
$scratch.Copy()
$ws.Range("E32").PasteSpecial(-4163)

$scratch.Value = "27.03"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)

$scratch.Value = "  +0.80%  "
$scratch.Copy()
$ws.Range("E33").PasteSpecial(-4163)

$scratch.Value = "  +0.55%  "
$scratch.Copy()
$ws.Range("E34").PasteSpecial(-4163)

$scratch.Value = "0.0₃0809"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)

$scratch.Value = "  +3.52%  "
$scratch.Copy()
$ws.Range("E35").PasteSpecial(-4163)

$scratch.Value = "  -0.94%  "
$scratch.Copy()
$ws.Range("E36").PasteSpecial(-4163)

$scratch.Value = "5.74"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)

$scratch.Value = "  +0.26%  "
$scratch.Copy()
$ws.Range("E37").PasteSpecial(-4163)

$scratch.Value = "50.17"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)

$scratch.Value = "  +0.39%  "
$scratch.Copy()
$ws.Range("E38").PasteSpecial(-4163)

$scratch.Value = "  -2.87%  "
$scratch.Copy()
$ws.Range("E39").PasteSpecial(-4163)

$scratch.Value = "8.98"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)

$scratch.Value = "  +0.41%  "
$scratch.Copy()
$ws.Range("E40").PasteSpecial(-4163)

$scratch.Value = "0.120"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)

$scratch.Value = "  +8.78%  "
$scratch.Copy()
$ws.Range("E41").PasteSpecial(-4163)

$scratch.Value = "2.84"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)

$scratch.Value = "  -2.50%  "
$scratch.Copy()
$ws.Range("E42").PasteSpecial(-4163)

$scratch.Value = "389.90"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)

$scratch.Value = "  -4.68%  "
$scratch.Copy()
$ws.Range("E43").PasteSpecial(-4163)

$scratch.Value = "39.12"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)

$scratch.Value = "  +1.63%  "
$scratch.Copy()
$ws.Range("E44").PasteSpecial(-4163)

$scratch.Value = "  -0.70%  "
$scratch.Copy()
$ws.Range("E45").PasteSpecial(-4163)

$scratch.Value = "  -4.25%  "
$scratch.Copy()
$ws.Range("E46").PasteSpecial(-4163)

$scratch.Value = "2.677.90"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)

$scratch.Value = "  -3.11%  "
$scratch.Copy()
$ws.Range("E47").PasteSpecial(-4163)

$scratch.Value = "130.12"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)

$scratch.Value = "  +2.23%  "
$scratch.Copy()
$ws.Range("E48").PasteSpecial(-4163)

$scratch.Value = "0.107"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)

$scratch.Value = "  -0.87%  "
$scratch.Copy()
$ws.Range("E50").PasteSpecial(-4163)

$scratch.Value = "  -0.61%  "
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)

# Remove the scratch cell entirely (contents + formatting) so it leaves
# no trace in the saved worksheet (used range / dimension unaffected).
$scratch.Clear()
$excel.CutCopyMode = $false
